$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 13: 45183 -> 45184
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
